$d = $word.ActiveDocument

# --- Hunk 1 ---------------------------------------------------------------
# The "Referente" bullet ends with a single underlined run whose text is
# " )" (a space followed by the closing parenthesis). Split it into two
# runs: the underlined " " stays as-is, and a brand-new, unformatted run
# carries the ")" character.
$r1 = $d.Content
$found1 = $r1.Find.Execute(" )", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $insertPoint = $d.Range($r1.End, $r1.End)
    $insertPoint.InsertAfter(")")
    $oldParen = $d.Range($r1.End - 1, $r1.End)
    $oldParen.Delete()
}

# --- Hunk 2 -----------------------------------------------------------------
# "Rilevazione" bullet: the three runs ", DataRil" / "ev" / ", DataIns,
# ModAcquisizione, InfoAmb" collapse into a single run with the combined
# (and shortened) text ", DataRil, DataIns, ModAcquisizione, InfoAmb".
$r2 = $d.Content
$r2.Find.Execute(", DataRilev, DataIns, ModAcquisizione, InfoAmb", $true, $false, $false, $false, $false, $true, 1, $false, ", DataRil, DataIns, ModAcquisizione, InfoAmb", 2)

# --- Hunk 3 -----------------------------------------------------------------
# Same bullet: ", RespRilev" becomes ", RespRil".
$r3 = $d.Content
$r3.Find.Execute(", RespRilev", $true, $false, $false, $false, $false, $true, 1, $false, ", RespRil", 2)
